$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to track several contacts (rows 2-5); now only a single
# (test) contact remains, so drop the old rows 3-5 entirely.
$ws.Range("A3:A5").EntireRow.Delete() | Out-Null

# Row 2: replace the contact with the new test data.
$ws.Range("A2").Value = "teste"
$ws.Range("B2").Value = "(61)98446-8993"

# The message text moves from E2 into D2 (kept verbatim), E2 is now just an
# empty cell, and F2 gets a standalone PDF path (no longer shares the old
# "PDF/BV NT - Catálogo.pdf" string).
$msg = "Olá Professor(a), tudo bem?`n`nAqui quem fala é Matheus Xavier da NT Editora.`n`nTomo a liberdade de encaminhar essa mensagem no intuito de compartilhar o nosso Catálogo de Livros e Catálogo descritivo da nossa Biblioteca Virtual para a Formação Técnica e Profissional.`n`nSabemos do desafio que é o fornecimento de materiais didáticos de qualidade e com foco no aluno do ensino profissional. Face ao exposto, gostaria de agendar uma rápida apresentação das nossas soluções educacionais que irão contribuir com a oferta da Vossa Instituição de Ensino.`n`nPodemos falar?"

$ws.Range("E2").ClearContents() | Out-Null

$ws.Range("D2").Value = $msg
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").WrapText = $true
$ws.Range("D2").RowHeight = 102.2

$ws.Range("F2").Value = "/home/gustavo/Desktop/meus_projetos/whatsapp_disparador/PDF/catalogo.pdf"

# Columns E/F grew noticeably wider to fit the longer message/path text.
$ws.Range("E1").ColumnWidth = 65.5
$ws.Range("F1").ColumnWidth = 65.5

$ws.Range("B2").Select() | Out-Null
